$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 868
$ws.Range("F3").Value = 991
$ws.Range("F4").Value = 777
$ws.Range("F6").Value = 437
$ws.Range("F7").Value = 676
$ws.Range("F8").Value = 153
$ws.Range("F9").Value = 1275
$ws.Range("F10").Value = 702
$ws.Range("F12").Value = 543
$ws.Range("F13").Value = 181
$ws.Range("F14").Value = 34
$ws.Range("F15").Value = 913
$ws.Range("F16").Value = 9
$ws.Range("F17").Value = 397
$ws.Range("F18").Value = 372
$ws.Range("F20").Value = 581
$ws.Range("F21").Value = 138
$ws.Range("F22").Value = 627
$ws.Range("F23").Value = 34
$ws.Range("F24").Value = 958
$ws.Range("F25").Value = 12

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 237
$ws.Range("F8").Value = 54
$ws.Range("F11").Value = 110

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 868
$ws.Range("F5").Value = 991
$ws.Range("F6").Value = 777
$ws.Range("F8").Value = 437
$ws.Range("F9").Value = 676
$ws.Range("F10").Value = 153
$ws.Range("F11").Value = 1275
$ws.Range("F12").Value = 702
$ws.Range("F16").Value = 543
$ws.Range("F18").Value = 181
$ws.Range("F19").Value = 34
$ws.Range("F20").Value = 913
$ws.Range("F22").Value = 9
$ws.Range("F23").Value = 397
$ws.Range("F24").Value = 372
$ws.Range("F26").Value = 237
$ws.Range("F27").Value = 54
$ws.Range("F28").Value = 581
$ws.Range("F31").Value = 110
$ws.Range("F32").Value = 110
$ws.Range("F33").Value = 138
$ws.Range("F34").Value = 627
$ws.Range("F35").Value = 34
$ws.Range("F36").Value = 958
$ws.Range("F37").Value = 12
